$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D:E").Insert(-4161, 0)

$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 820200
$ws.Range("E8").Value = 466900
$ws.Range("D9").Value = 609400
$ws.Range("E9").Value = 448500
$ws.Range("D10").Value = 210800
$ws.Range("E10").Value = 18400
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 45700
$ws.Range("E15").Value = 46800
$ws.Range("D17").Value = 732100
$ws.Range("E17").Value = 456500
$ws.Range("D18").Value = 88100
$ws.Range("E18").Value = 10400
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 133800
$ws.Range("E21").Value = 57200
$ws.Range("D22").Value = 42400
$ws.Range("E22").Value = 41100
$ws.Range("D23").Value = 45700
$ws.Range("E23").Value = -30700
$ws.Range("D24").Value = 400
$ws.Range("E24").Value = 600
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 45300
$ws.Range("E26").Value = -31300
$ws.Range("D27").Value = 32700
$ws.Range("E27").Value = -42300
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 32700
$ws.Range("E33").Value = -42300
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 32700
$ws.Range("E35").Value = -42300
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 8300
$ws.Range("E41").Value = 6900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 347200
$ws.Range("E43").Value = 209800
$ws.Range("D44").Value = 128900
$ws.Range("E44").Value = 130500
$ws.Range("D45").Value = 66100
$ws.Range("E45").Value = 104700
$ws.Range("D46").Value = 550500
$ws.Range("E46").Value = 451900
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 1141100
$ws.Range("E48").Value = 1148400
$ws.Range("D49").Value = 2273100
$ws.Range("E49").Value = 2283300
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 40700
$ws.Range("E52").Value = 42300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 4005500
$ws.Range("E54").Value = 3925800
$ws.Range("D57").Value = 171800
$ws.Range("E57").Value = 137100
$ws.Range("D58").Value = 376800
$ws.Range("E58").Value = 240600
$ws.Range("D59").Value = 265200
$ws.Range("E59").Value = 300400
$ws.Range("D60").Value = 813900
$ws.Range("E60").Value = 678000
$ws.Range("D61").Value = 2560700
$ws.Range("E61").Value = 2561000
$ws.Range("D62").Value = 118100
$ws.Range("E62").Value = 117100
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 3524900
$ws.Range("E66").Value = 3389200
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 12100
$ws.Range("E70").Value = 12700
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("E72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 468500
$ws.Range("E76").Value = 523900
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 32700
$ws.Range("E81").Value = -42300
$ws.Range("D83").Value = 45700
$ws.Range("E83").Value = 46800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -3900
$ws.Range("E89").Value = 65900
$ws.Range("D91").Value = -31000
$ws.Range("E91").Value = -28400
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -27400
$ws.Range("E94").Value = -12200
$ws.Range("D96").Value = -100700
$ws.Range("E96").Value = -100700
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 32700
$ws.Range("E100").Value = -52000
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 1400
$ws.Range("E102").Value = 1800
